# tidsregistrering for d. 20/03
# Append the time-registration entries for 17/03 (continuation) and 20/03
# below the existing table (rows 38-47), widen column F for the longer
# activity text, and move the view/selection down to where the new rows
# were entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new data rows -----------------------------------------------------
# Row => Date (serial, $null = none), Role (E, $null = none), F, G, H
$rows = @(
    @{ Row = 38; Date = 42811; Role = "Implementer"; F = "Ændre i brugergrænseFladen";                 G = "8.20";  H = "11.20" }
    @{ Row = 39; Date = 42814; Role = $null;          F = "OC 10 exporter pdf";                        G = "8.25";  H = "9.05" }
    @{ Row = 40; Date = $null; Role = $null;          F = "DOM UC 9";                                  G = "8.25";  H = "9.05" }
    @{ Row = 41; Date = $null; Role = $null;          F = "DOM samlet";                                G = "8.50";  H = "9.05" }
    @{ Row = 42; Date = $null; Role = $null;          F = "Review OC 13";                              G = "9.30";  H = "9.55" }
    @{ Row = 43; Date = $null; Role = $null;          F = "Test suite OC 13 beregnBøjningsMoment";     G = "10.05"; H = "10.55" }
    @{ Row = 44; Date = $null; Role = $null;          F = "Dataorbog samling";                         G = "11.00"; H = "11.30" }
    @{ Row = 45; Date = $null; Role = $null;          F = "Opgaver til Brugertest";                    G = "12.25"; H = "12.55" }
    @{ Row = 46; Date = $null; Role = $null;          F = "Review på OC 8 test suite";                 G = "12.05"; H = "13.15" }
    @{ Row = 47; Date = $null; Role = $null;          F = "Review på OC 11 DCD, SD";                   G = "13.15"; H = "14.30" }
)

# G:H hold clock-times typed as plain text (same convention already used
# throughout this sheet, e.g. "8.30", "15.30" in row 4). Flip the whole
# target block to text first so the "8.20"-style values aren't
# auto-parsed as numbers, then restore the default ("Normal") style so no
# lingering number format is left on the cells.
$timeBlock = $ws.Range("G38:H47")
$timeBlock.NumberFormat = "@"

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($null -ne $r.Date) {
        $dateCell = $ws.Cells.Item($rowNum, 1)
        $dateCell.NumberFormat = "m/d/yy"
        $dateCell.Value = $r.Date
    }

    if ($null -ne $r.Role) {
        $ws.Cells.Item($rowNum, 5).Value = $r.Role
    }

    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
}

$timeBlock.Style = "Normal"

# --- widen column F for the longer activity text -----------------------
$ws.Columns.Item(6).ColumnWidth = 37.6

# --- move the view to follow the newly entered rows ---------------------
$win = $excel.ActiveWindow
$ws.Activate()
$win.ScrollRow = 22
$win.ScrollColumn = 4
$ws.Range("H48").Select()
